$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 61
$ws.Range("A61").Value = 60.0
$ws.Range("B61").Value = "Tuesday, Jan 10"
$ws.Range("C61").Value = "4:30 AM"
$ws.Range("D61").Value = "UNKNOWN"
$ws.Range("E61").Value = "Warsaw"
$ws.Range("F61").Value = "(WAW)"
$ws.Range("G61").Value = "Enter Air "
$ws.Range("H61").Value = "B738"
$ws.Range("I61").Value = "(SP-ESB)"
$ws.Range("J61").Value = "4:56 AM"
$ws.Range("L61").Value = "0 hours, 26 minutes"

# Row 62
$ws.Range("A62").Value = 61.0
$ws.Range("B62").Value = "Tuesday, Jan 10"
$ws.Range("C62").Value = "5:40 AM"
$ws.Range("D62").Value = "LO3942"
$ws.Range("E62").Value = "Warsaw"
$ws.Range("F62").Value = "(WAW)"
$ws.Range("G62").Value = "LOT "
$ws.Range("H62").Value = "E170"
$ws.Range("I62").Value = "(SP-LDF)"
$ws.Range("J62").Value = "5:46 AM"
$ws.Range("L62").Value = "0 hours, 6 minutes"

# Row 63
$ws.Range("A63").Value = 62.0
$ws.Range("B63").Value = "Tuesday, Jan 10"
$ws.Range("C63").Value = "5:50 AM"
$ws.Range("D63").Value = "FR7895"
$ws.Range("E63").Value = "Paris"
$ws.Range("F63").Value = "(BVA)"
$ws.Range("G63").Value = "Ryanair "
$ws.Range("H63").Value = "B738"
$ws.Range("I63").Value = "(SP-RSM)"
$ws.Range("J63").Value = "6:04 AM"
$ws.Range("L63").Value = "0 hours, 14 minutes"

# Row 64
$ws.Range("A64").Value = 63.0
$ws.Range("B64").Value = "Tuesday, Jan 10"
$ws.Range("C64").Value = "6:30 AM"
$ws.Range("D64").Value = "LH1381"
$ws.Range("E64").Value = "Frankfurt"
$ws.Range("F64").Value = "(FRA)"
$ws.Range("G64").Value = "Lufthansa "
$ws.Range("H64").Value = "CRJ9"
$ws.Range("I64").Value = "(D-ACNT)"
$ws.Range("J64").Value = "6:38 AM"
$ws.Range("L64").Value = "0 hours, 8 minutes"

# Row 65
$ws.Range("A65").Value = 64.0
$ws.Range("B65").Value = "Tuesday, Jan 10"
$ws.Range("C65").Value = "8:05 AM"
$ws.Range("D65").Value = "LH1641"
$ws.Range("E65").Value = "Munich"
$ws.Range("F65").Value = "(MUC)"
$ws.Range("G65").Value = "Lufthansa "
$ws.Range("H65").Value = "CRJ9"
$ws.Range("I65").Value = "(D-ACKH)"
$ws.Range("J65").Value = "8:11 AM"
$ws.Range("L65").Value = "0 hours, 6 minutes"

# Row 66
$ws.Range("A66").Value = 65.0
$ws.Range("B66").Value = "Tuesday, Jan 10"
$ws.Range("C66").Value = "8:30 AM"
$ws.Range("D66").Value = "FR5000"
$ws.Range("E66").Value = "Brussels"
$ws.Range("F66").Value = "(CRL)"
$ws.Range("G66").Value = "Ryanair "
$ws.Range("H66").Value = "B738"
$ws.Range("I66").Value = "(EI-DPG)"
$ws.Range("J66").Value = "8:55 AM"
$ws.Range("L66").Value = "0 hours, 25 minutes"

# Row 67
$ws.Range("A67").Value = 66.0
$ws.Range("B67").Value = "Tuesday, Jan 10"
$ws.Range("C67").Value = "10:50 AM"
$ws.Range("D67").Value = "FR3728"
$ws.Range("E67").Value = "Billund"
$ws.Range("F67").Value = "(BLL)"
$ws.Range("G67").Value = "Ryanair "
$ws.Range("H67").Value = "B738"
$ws.Range("I67").Value = "(SP-RSM)"
$ws.Range("J67").Value = "10:47 AM"
$ws.Range("L67").Value = "0 hours, -3 minutes"

# Row 68
$ws.Range("A68").Value = 67.0
$ws.Range("B68").Value = "Tuesday, Jan 10"
$ws.Range("C68").Value = "11:10 AM"
$ws.Range("D68").Value = "UNKNOWN"
$ws.Range("E68").Value = "Poprad"
$ws.Range("F68").Value = "(TAT)"
$ws.Range("G68").Value = "AMC Aviation "
$ws.Range("H68").Value = "PC24"
$ws.Range("I68").Value = "(SP-AGA)"
$ws.Range("J68").Value = "12:42 PM"
$ws.Range("L68").Value = "1 hours, 32 minutes"

# Row 69
$ws.Range("A69").Value = 68.0
$ws.Range("B69").Value = "Tuesday, Jan 10"
$ws.Range("C69").Value = "11:50 AM"
$ws.Range("D69").Value = "LO3946"
$ws.Range("E69").Value = "Warsaw"
$ws.Range("F69").Value = "(WAW)"
$ws.Range("G69").Value = "LOT "
$ws.Range("H69").Value = "E75S"
$ws.Range("I69").Value = "(SP-LIB)"
$ws.Range("J69").Value = "12:16 PM"
$ws.Range("L69").Value = "0 hours, 26 minutes"

# Row 70
$ws.Range("A70").Value = 69.0
$ws.Range("B70").Value = "Tuesday, Jan 10"
$ws.Range("C70").Value = "12:55 PM"
$ws.Range("D70").Value = "LH1637"
$ws.Range("E70").Value = "Munich"
$ws.Range("F70").Value = "(MUC)"
$ws.Range("G70").Value = "Lufthansa "
$ws.Range("H70").Value = "CRJ9"
$ws.Range("I70").Value = "(D-ACNN)"
$ws.Range("J70").Value = "1:00 PM"
$ws.Range("L70").Value = "0 hours, 5 minutes"

# Row 71
$ws.Range("A71").Value = 70.0
$ws.Range("B71").Value = "Tuesday, Jan 10"
$ws.Range("C71").Value = "1:20 PM"
$ws.Range("D71").Value = "LO6531"
$ws.Range("E71").Value = "Puerto Plata"
$ws.Range("F71").Value = "(POP)"
$ws.Range("G71").Value = "LOT "
$ws.Range("H71").Value = "B788"
$ws.Range("I71").Value = "(SP-LRD)"
$ws.Range("J71").Value = "1:20 PM"
$ws.Range("L71").Value = "0 hours, 0 minutes"

# Row 72
$ws.Range("A72").Value = 71.0
$ws.Range("B72").Value = "Tuesday, Jan 10"
$ws.Range("C72").Value = "2:35 PM"
$ws.Range("D72").Value = "FR7898"
$ws.Range("E72").Value = "Amman"
$ws.Range("F72").Value = "(AMM)"
$ws.Range("G72").Value = "Ryanair "
$ws.Range("H72").Value = "B738"
$ws.Range("I72").Value = "(SP-RKR)"
$ws.Range("J72").Value = "2:42 PM"
$ws.Range("L72").Value = "0 hours, 7 minutes"

# Row 73
$ws.Range("A73").Value = 72.0
$ws.Range("B73").Value = "Tuesday, Jan 10"
$ws.Range("C73").Value = "2:45 PM"
$ws.Range("D73").Value = "LO3944"
$ws.Range("E73").Value = "Warsaw"
$ws.Range("F73").Value = "(WAW)"
$ws.Range("G73").Value = "LOT "
$ws.Range("H73").Value = "E75S"
$ws.Range("I73").Value = "(SP-LIA)"
$ws.Range("J73").Value = "2:53 PM"
$ws.Range("L73").Value = "0 hours, 8 minutes"

# Row 74
$ws.Range("A74").Value = 73.0
$ws.Range("B74").Value = "Tuesday, Jan 10"
$ws.Range("C74").Value = "3:25 PM"
$ws.Range("D74").Value = "LH1391"
$ws.Range("E74").Value = "Frankfurt"
$ws.Range("F74").Value = "(FRA)"
$ws.Range("G74").Value = "Lufthansa "
$ws.Range("H74").Value = "CRJ9"
$ws.Range("I74").Value = "(D-ACNJ)"
$ws.Range("J74").Value = "3:36 PM"
$ws.Range("L74").Value = "0 hours, 11 minutes"

# Row 75
$ws.Range("A75").Value = 74.0
$ws.Range("B75").Value = "Tuesday, Jan 10"
$ws.Range("C75").Value = "6:20 PM"
$ws.Range("D75").Value = "LO3948"
$ws.Range("E75").Value = "Warsaw"
$ws.Range("F75").Value = "(WAW)"
$ws.Range("G75").Value = "LOT "
$ws.Range("H75").Value = "E170"
$ws.Range("I75").Value = "(SP-LDF)"
$ws.Range("J75").Value = "6:20 PM"
$ws.Range("L75").Value = "0 hours, 0 minutes"

# Row 76
$ws.Range("A76").Value = 75.0
$ws.Range("B76").Value = "Tuesday, Jan 10"
$ws.Range("C76").Value = "9:30 PM"
$ws.Range("D76").Value = "FR1751"
$ws.Range("E76").Value = "London"
$ws.Range("F76").Value = "(STN)"
$ws.Range("G76").Value = "Ryanair "
$ws.Range("H76").Value = "B738"
$ws.Range("I76").Value = "(EI-EMJ)"
$ws.Range("J76").Value = "9:30 PM"
$ws.Range("L76").Value = "0 hours, 0 minutes"

# Row 77
$ws.Range("A77").Value = 76.0
$ws.Range("B77").Value = "Tuesday, Jan 10"
$ws.Range("C77").Value = "9:40 PM"
$ws.Range("D77").Value = "P81987"
$ws.Range("E77").Value = "Cologne"
$ws.Range("F77").Value = "(CGN)"
$ws.Range("G77").Value = "SprintAir "
$ws.Range("H77").Value = "AT72"
$ws.Range("I77").Value = "(SP-SPG)"
$ws.Range("J77").Value = "9:45 PM"
$ws.Range("L77").Value = "0 hours, 5 minutes"
